$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 updates (shift B..E content in: B26 keeps 8:30am already there,
# but per diff new F26 gets a value, G26 gets status, I26 becomes numeric 34,
# M26 gets a new note)
$ws.Range("F26").Value = "6:31pm"
$ws.Range("G26").Value = "complete"
$ws.Range("I26").Value = 34
$ws.Range("M26").Value = "12 days"

# Row 27 new entries
$ws.Range("A27").Value = 45559
$ws.Range("B27").Value = "5:36am"
$ws.Range("D27").Value = "6:07pm"
$ws.Range("E27").Value = "6:28pm"
$ws.Range("O27").Value = "SIT FOR 10 HOUS AND STUDY---BOYCOTT EVERYTHING"

# Row 28 updates
$ws.Range("A28").Value = "25th-Sep"
$ws.Range("D28").Value = "7:14AM"
$ws.Range("E28").Value = "8:07am"
$ws.Range("F28").Value = "8:47am"
$ws.Range("I28").Value = 29
$ws.Range("M28").Value = "10 days"
$ws.Range("N28").Value = "1 greedy , 6 interval pending"

# Update the active selection to I26 (was I34)
$ws.Range("I26").Select()
